$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("B2").Value = 0.006876353814593728
$ws.Range("C2").Value = 0.000002220651329265522
$ws.Range("D2").Value = 157.8057217802531
$ws.Range("E2").Value = 6.48142807727062
$ws.Range("G2").Value = 164.2940284319897

# Row 3 updates
$ws.Range("B3").Value = 3.182878228561681
$ws.Range("C3").Value = 1.65323645889881
$ws.Range("D3").Value = 0.1529057820181812
$ws.Range("E3").Value = 0.4998867070740569
$ws.Range("G3").Value = 5.488907176552729
